{"js": "// Replace the trailing empty paragraph (the last paragraph in the main\n// body, immediately before the section break) with a paragraph made of\n// two runs:\n//   1) \"Version management\"\n//   2) \" is a way to manage and track changes related to a collection\n//      of entities.\"\n// Using Range.insertOoxml(..., \"Replace\") lets us control the exact run\n// split (matching the authored OOXML) instead of letting consecutive\n// insertText() calls coalesce into a single run.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst items = paragraphs.items;\nconst target = items[items.length - 1];\nconst targetRange = target.getRange();\n\nconst flatOpcXml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" ' +\n  'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p>' +\n  '<w:r><w:t>Version management</w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\"> is a way to manage and track changes related to a collection of entities.</w:t></w:r>' +\n  '</w:p>' +\n  '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>';\n\ntargetRange.insertOoxml(flatOpcXml, \"Replace\");\nawait context.sync();\n", "ps1": "# Replace the trailing empty paragraph (the last paragraph in the main\n# body, immediately before the section break) with a paragraph made of\n# two runs:\n#   1) \"Version management\"\n#   2) \" is a way to manage and track changes related to a collection\n#      of entities.\"\n# Range.InsertXML(xml, \"Replace\") lets us control the exact run split\n# (matching the authored OOXML) instead of letting InsertAfter calls\n# coalesce into a single run.\n\n$d = $word.ActiveDocument\n\n$count = $d.Paragraphs.Count\n$target = $d.Paragraphs.Item($count)\n$r = $target.Range\n\n$flatOpcXml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' + `\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' + `\n  '<pkg:part pkg:name=\"/word/document.xml\" ' + `\n  'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' + `\n  '<pkg:xmlData>' + `\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' + `\n  '<w:body>' + `\n  '<w:p>' + `\n  '<w:r><w:t>Version management</w:t></w:r>' + `\n  '<w:r><w:t xml:space=\"preserve\"> is a way to manage and track changes related to a collection of entities.</w:t></w:r>' + `\n  '</w:p>' + `\n  '</w:body>' + `\n  '</w:document>' + `\n  '</pkg:xmlData></pkg:part></pkg:package>'\n\n$r.InsertXML($flatOpcXml, \"Replace\")\n"}
